{"js": "// This script updates the two-digit multiplication problems in the\n// document's table cells, replacing each old \"AA\u00d7BB=\" expression with\n// its new value, per the commit's regenerated worksheet numbers.\n\nconst replacements = [\n  [\"93\u00d749=\", \"12\u00d774=\"],\n  [\"93\u00d717=\", \"49\u00d747=\"],\n  [\"70\u00d716=\", \"88\u00d794=\"],\n  [\"45\u00d756=\", \"88\u00d723=\"],\n  [\"12\u00d792=\", \"51\u00d735=\"],\n  [\"85\u00d770=\", \"38\u00d752=\"],\n  [\"20\u00d798=\", \"75\u00d792=\"],\n  [\"73\u00d762=\", \"49\u00d757=\"],\n  [\"33\u00d722=\", \"79\u00d740=\"],\n  [\"68\u00d747=\", \"60\u00d796=\"],\n  [\"19\u00d746=\", \"58\u00d744=\"],\n  [\"23\u00d772=\", \"45\u00d784=\"],\n  [\"63\u00d739=\", \"60\u00d754=\"],\n  [\"75\u00d749=\", \"85\u00d774=\"],\n  [\"76\u00d724=\", \"33\u00d780=\"],\n  [\"91\u00d716=\", \"19\u00d733=\"],\n  [\"51\u00d739=\", \"16\u00d763=\"],\n  [\"40\u00d737=\", \"70\u00d756=\"],\n  [\"15\u00d748=\", \"89\u00d739=\"],\n  [\"60\u00d725=\", \"64\u00d792=\"],\n  [\"78\u00d761=\", \"98\u00d718=\"],\n  [\"69\u00d716=\", \"20\u00d716=\"],\n  [\"60\u00d794=\", \"19\u00d762=\"],\n  [\"69\u00d754=\", \"91\u00d775=\"],\n  [\"79\u00d752=\", \"84\u00d777=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# This script updates the two-digit multiplication problems in the\n# document's table cells, replacing each old \"AA\u00d7BB=\" expression with\n# its new value, per the commit's regenerated worksheet numbers.\n\n$replacements = @(\n    @(\"93\u00d749=\", \"12\u00d774=\"),\n    @(\"93\u00d717=\", \"49\u00d747=\"),\n    @(\"70\u00d716=\", \"88\u00d794=\"),\n    @(\"45\u00d756=\", \"88\u00d723=\"),\n    @(\"12\u00d792=\", \"51\u00d735=\"),\n    @(\"85\u00d770=\", \"38\u00d752=\"),\n    @(\"20\u00d798=\", \"75\u00d792=\"),\n    @(\"73\u00d762=\", \"49\u00d757=\"),\n    @(\"33\u00d722=\", \"79\u00d740=\"),\n    @(\"68\u00d747=\", \"60\u00d796=\"),\n    @(\"19\u00d746=\", \"58\u00d744=\"),\n    @(\"23\u00d772=\", \"45\u00d784=\"),\n    @(\"63\u00d739=\", \"60\u00d754=\"),\n    @(\"75\u00d749=\", \"85\u00d774=\"),\n    @(\"76\u00d724=\", \"33\u00d780=\"),\n    @(\"91\u00d716=\", \"19\u00d733=\"),\n    @(\"51\u00d739=\", \"16\u00d763=\"),\n    @(\"40\u00d737=\", \"70\u00d756=\"),\n    @(\"15\u00d748=\", \"89\u00d739=\"),\n    @(\"60\u00d725=\", \"64\u00d792=\"),\n    @(\"78\u00d761=\", \"98\u00d718=\"),\n    @(\"69\u00d716=\", \"20\u00d716=\"),\n    @(\"60\u00d794=\", \"19\u00d762=\"),\n    @(\"69\u00d754=\", \"91\u00d775=\"),\n    @(\"79\u00d752=\", \"84\u00d777=\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
